# cropping coordinates calculate for dataset05, with my py pandas script
#
# Rename several header labels in row 1 of the cropping table:
#   - drop the "(coordinate)" suffix from the z index/interest labels
#   - rename the bounding-box size columns from X/Y/Z (width/height/depth)
#     to the x_size / y_size / z_size naming convention

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AE1").Value = "z_interest"
$ws.Range("AB1").Value = "z_size"
$ws.Range("AA1").Value = "y_size"
$ws.Range("Z1").Value  = "x_size"
$ws.Range("X1").Value  = "z1_ind"
$ws.Range("W1").Value  = "z0_ind"
